# Fruta / hortaliza, semanal
#
# Two new price records (Terminal Hortofrutícola Agro Chillán - Pera,
# Packham's Triumph, date 44495) are inserted into the daily data table,
# right after the existing row for date 44424 (worksheet row 58). This
# pushes every subsequent record down by two rows and grows the sheet
# from 148 to 150 data-bearing rows (dimension A1:T148 -> A1:T150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before row 59, shifting rows 59:148 down to 61:150.
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

function Set-DataRow {
    param($r, $vals)

    $ws.Cells.Item($r, 1).Value2  = $vals[0]   # Mercado ID
    $ws.Cells.Item($r, 2).Value2  = $vals[1]   # Mercado
    $ws.Cells.Item($r, 3).Value2  = $vals[2]   # Región
    $ws.Cells.Item($r, 4).Value2  = $vals[3]   # Fecha
    $ws.Cells.Item($r, 5).Value2  = $vals[4]   # Codreg
    $ws.Cells.Item($r, 6).Value2  = $vals[5]   # Tipo
    $ws.Cells.Item($r, 7).Value2  = $vals[6]   # Producto ID
    $ws.Cells.Item($r, 8).Value2  = $vals[7]   # Producto
    $ws.Cells.Item($r, 9).Value2  = $vals[8]   # Categoría ID
    $ws.Cells.Item($r, 10).Value2 = $vals[9]   # Categoría
    $ws.Cells.Item($r, 11).Value2 = $vals[10]  # Variedad
    $ws.Cells.Item($r, 12).Value2 = $vals[11]  # Calidad
    $ws.Cells.Item($r, 13).Value2 = $vals[12]  # Volumen
    $ws.Cells.Item($r, 14).Value2 = $vals[13]  # Precio mínimo
    $ws.Cells.Item($r, 15).Value2 = $vals[14]  # Precio máximo
    $ws.Cells.Item($r, 16).Value2 = $vals[15]  # Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value2 = $vals[16]  # Unidad de comercialización
    $ws.Cells.Item($r, 18).Value2 = $vals[17]  # Origen
    $ws.Cells.Item($r, 19).Value2 = $vals[18]  # Precio $/Kg
    $ws.Cells.Item($r, 20).Value2 = $vals[19]  # Kg / unidad
}

# New record: calidad "Especial"
Set-DataRow 59 @(
    7, "Terminal Hortofrutícola Agro Chillán", "Ñuble",
    44495, 16, "Fruta", 100104, "Frutos de pepita", 100104005, "Pera",
    "Packham's Triumph", "Especial", 80, 10000, 10000, 10000,
    "`$/caja 16 kilos empedrada", "Provincia de Curicó", 625, 16
)

# New record: calidad "Primera"
Set-DataRow 60 @(
    7, "Terminal Hortofrutícola Agro Chillán", "Ñuble",
    44495, 16, "Fruta", 100104, "Frutos de pepita", 100104005, "Pera",
    "Packham's Triumph", "Primera", 160, 8800, 9000, 8900,
    "`$/caja 16 kilos empedrada", "Provincia de Curicó", 556, 16
)

$addr = $ws.UsedRange.Address()
"Final used range: $addr"
